$d = $word.ActiveDocument
$d.Content.Find.Execute("subsidiary, subsidiaries, affiliate, acquiror, merger sub, covenantor, acquired company, acquiring company, surviving corporation, surviving company", $true, $false, $false, $false, $false,
                         $true, 1, $false, "subsidiary, subsidiaries, affiliate, acquiror, merger sub, covenantor, acquired company, acquiring company, surviving corporation, surviving company", 2)
